$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1457.8235
$ws.Range("I100").Value = 1253
$ws.Range("K100").Value = 1253
$ws.Range("M100").Value = -712
$ws.Range("H108").Value = 31618.666
$ws.Range("J108").Value = 31618.666
$ws.Range("L108").Value = 31618.666
$ws.Range("N108").Value = -39298.666
$ws.Range("H123").Value = 37250
$ws.Range("J123").Value = 37250
$ws.Range("L123").Value = 37250
$ws.Range("N123").Value = -47050
$ws.Range("H126").Value = 46766.668
$ws.Range("J126").Value = 46766.668
$ws.Range("L126").Value = 46766.668
$ws.Range("N126").Value = -56646.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 77504.5
$ws.Range("I23").Value = 76672.664
$ws.Range("J23").Value = 80000
$ws.Range("K23").Value = 76672.664
$ws.Range("L23").Value = 80000
$ws.Range("M23").Value = -76413.664
$ws.Range("N23").Value = -80518
$ws.Range("H111").Value = 49640
$ws.Range("J111").Value = 49640
$ws.Range("L111").Value = 49640
$ws.Range("N111").Value = -57820
$ws.Range("H137").Value = 31127.25
$ws.Range("J137").Value = 41266.668
$ws.Range("L137").Value = 41266.668
$ws.Range("N137").Value = -51466.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 43533.332
$ws.Range("J59").Value = 43533.332
$ws.Range("L59").Value = 43533.332
$ws.Range("N59").Value = -45227.332
$ws.Range("H108").Value = 47644
$ws.Range("J108").Value = 47644
$ws.Range("L108").Value = 47644
$ws.Range("N108").Value = -55324
$ws.Range("H110").Value = 46658.332
$ws.Range("J110").Value = 46658.332
$ws.Range("L110").Value = 46658.332
$ws.Range("N110").Value = -54838.332
$ws.Range("H111").Value = 47702
$ws.Range("J111").Value = 47702
$ws.Range("L111").Value = 47702
$ws.Range("N111").Value = -55882
$ws.Range("H112").Value = 45996
$ws.Range("J112").Value = 45996
$ws.Range("L112").Value = 45996
$ws.Range("N112").Value = -48950
$ws.Range("H116").Value = 43459.5
$ws.Range("J116").Value = 43459.5
$ws.Range("L116").Value = 43459.5
$ws.Range("N116").Value = -52637.5
$ws.Range("H133").Value = 54999.332
$ws.Range("J133").Value = 54999.332
$ws.Range("L133").Value = 54999.332
$ws.Range("N133").Value = -65119.332
$ws.Range("H139").Value = 49500
$ws.Range("J139").Value = 49500
$ws.Range("L139").Value = 49500
$ws.Range("N139").Value = -59780

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 42254
$ws.Range("J118").Value = 42254
$ws.Range("L118").Value = 42254
$ws.Range("N118").Value = -45568
$ws.Range("H132").Value = 40527.082
$ws.Range("I132").Value = 1919.8462
$ws.Range("J132").Value = 131780.55
$ws.Range("K132").Value = 5759.5386
$ws.Range("L132").Value = 395341.65
$ws.Range("M132").Value = -3229.5386
$ws.Range("N132").Value = -400401.65
$ws.Range("H133").Value = 26474
$ws.Range("J133").Value = 26474
$ws.Range("L133").Value = 26474
$ws.Range("N133").Value = -31534
$ws.Range("H137").Value = 20058
$ws.Range("J137").Value = 20058
$ws.Range("L137").Value = 20058
$ws.Range("N137").Value = -30258
$ws.Range("H139").Value = 35500
$ws.Range("J139").Value = 29000
$ws.Range("L139").Value = 29000
$ws.Range("N139").Value = -39280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3632.425
$ws.Range("I5").Value = 13104.125
$ws.Range("J5").Value = 1264.5
$ws.Range("K5").Value = 39312.375
$ws.Range("L5").Value = 3793.5
$ws.Range("M5").Value = -39200.375
$ws.Range("N5").Value = -4017.5
$ws.Range("H101").Value = 5000
$ws.Range("J101").Value = 5000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -19868
$ws.Range("H131").Value = 2655.603
$ws.Range("I131").Value = 6623.25
$ws.Range("J131").Value = 1434.7885
$ws.Range("K131").Value = 19869.75
$ws.Range("L131").Value = 4304.3655
$ws.Range("M131").Value = -14829.75
$ws.Range("N131").Value = -14384.3655
$ws.Range("H132").Value = 2584.48
$ws.Range("J132").Value = 3041.3333
$ws.Range("L132").Value = 27371.9997
$ws.Range("N132").Value = -32431.9997
$ws.Range("H135").Value = 3632.425
$ws.Range("I135").Value = 13104.125
$ws.Range("J135").Value = 1264.5
$ws.Range("K135").Value = 117937.125
$ws.Range("L135").Value = 11380.5
$ws.Range("M135").Value = -115402.125
$ws.Range("N135").Value = -16450.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 47753
$ws.Range("J119").Value = 47753
$ws.Range("L119").Value = 47753
$ws.Range("N119").Value = -57429
$ws.Range("H124").Value = 39780
$ws.Range("J124").Value = 39780
$ws.Range("L124").Value = 39780
$ws.Range("N124").Value = -49600
$ws.Range("H128").Value = 37363
$ws.Range("J128").Value = 37363
$ws.Range("L128").Value = 37363
$ws.Range("N128").Value = -47323
$ws.Range("H137").Value = 19207.691
$ws.Range("J137").Value = 19207.691
$ws.Range("L137").Value = 19207.691
$ws.Range("N137").Value = -29407.691

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1360.1
$ws.Range("I22").Value = 800.1667
$ws.Range("J22").Value = 2200
$ws.Range("K22").Value = 800.1667
$ws.Range("L22").Value = 2200
$ws.Range("M22").Value = -505.1667
$ws.Range("N22").Value = -2790
$ws.Range("H27").Value = 1360.1
$ws.Range("I27").Value = 800.1667
$ws.Range("J27").Value = 2200
$ws.Range("K27").Value = 800.1667
$ws.Range("L27").Value = 2200
$ws.Range("M27").Value = -693.1667
$ws.Range("N27").Value = -2414
$ws.Range("H93").Value = 29414030
$ws.Range("I93").Value = 83334780
$ws.Range("J93").Value = 2708.818
$ws.Range("K93").Value = 83334780
$ws.Range("L93").Value = 2708.818
$ws.Range("M93").Value = -83333532
$ws.Range("N93").Value = -5204.818
$ws.Range("H100").Value = 2231.2104
$ws.Range("I100").Value = 1914.8462
$ws.Range("K100").Value = 1914.8462
$ws.Range("M100").Value = -1373.8462
$ws.Range("H105").Value = 43681.5
$ws.Range("J105").Value = 43681.5
$ws.Range("L105").Value = 43681.5
$ws.Range("N105").Value = -50669.5
$ws.Range("H116").Value = 47676
$ws.Range("J116").Value = 47676
$ws.Range("L116").Value = 47676
$ws.Range("N116").Value = -56854
$ws.Range("H120").Value = 51741.5
$ws.Range("J120").Value = 51741.5
$ws.Range("L120").Value = 51741.5
$ws.Range("N120").Value = -61417.5
$ws.Range("H132").Value = 4676.615
$ws.Range("I132").Value = 3059.8
$ws.Range("K132").Value = 9179.400000000001
$ws.Range("M132").Value = -6649.400000000001
$ws.Range("H133").Value = 21271.062
$ws.Range("J133").Value = 21271.062
$ws.Range("L133").Value = 21271.062
$ws.Range("N133").Value = -26331.062
$ws.Range("H137").Value = 23279.166
$ws.Range("J137").Value = 23279.166
$ws.Range("L137").Value = 23279.166
$ws.Range("N137").Value = -33479.166
$ws.Range("H139").Value = 42226.582
$ws.Range("J139").Value = 42226.582
$ws.Range("L139").Value = 42226.582
$ws.Range("N139").Value = -52506.582

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 32377
$ws.Range("J109").Value = 32377
$ws.Range("L109").Value = 32377
$ws.Range("N109").Value = -35151
$ws.Range("H117").Value = 47206
$ws.Range("J117").Value = 47206
$ws.Range("L117").Value = 47206
$ws.Range("N117").Value = -56384
$ws.Range("H139").Value = 19031.924
$ws.Range("J139").Value = 19031.924
$ws.Range("L139").Value = 19031.924
$ws.Range("N139").Value = -29311.924
